$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix casing: "designationCategoryDef" -> "DesignationCategoryDef"
$ws.Range("A26").Value = "DesignationCategoryDef+Decor.label"
$ws.Range("A27").Value = "DesignationCategoryDef+Decor.description"
$ws.Range("B26").Value = "DesignationCategoryDef"
$ws.Range("B27").Value = "DesignationCategoryDef"

# Update the active selection to match the saved view state
$ws.Range("E31").Select()
